$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 42.73245000928255
$ws.Range("B3").Value = 64.48049027741986
$ws.Range("B4").Value = 71.36781907553818
$ws.Range("H5").Value = 96.05163328668019
$ws.Range("H6").Value = 96.0505924644102
$ws.Range("H7").Value = 96.03847318460804
$ws.Range("C8").Value = 95.78703591613193
$ws.Range("C9").Value = 93.87418676168402
$ws.Range("C10").Value = 95.63664062851603
$ws.Range("D11").Value = 99.35103308912282
$ws.Range("D12").Value = 99.30725849218793
$ws.Range("D13").Value = 99.45435956309598
$ws.Range("E14").Value = 98.89222711877861
$ws.Range("E15").Value = 98.96192624723412
$ws.Range("E16").Value = 98.90203522729797
$ws.Range("F17").Value = 98.4021473454396
$ws.Range("F18").Value = 98.47376201889568
$ws.Range("F19").Value = 98.37427544606669
$ws.Range("G20").Value = 97.43640106241376
$ws.Range("G21").Value = 97.51190816632631
$ws.Range("G22").Value = 97.47220867545362
$ws.Range("B23").Value = 70.80629415240057
$ws.Range("B24").Value = 77.67162260235956
$ws.Range("H25").Value = 96.08848131486407
$ws.Range("H26").Value = 96.04496534387852
$ws.Range("C27").Value = 95.39568080142431
$ws.Range("C28").Value = 95.17044764034401
$ws.Range("D29").Value = 99.37089676117409
$ws.Range("D30").Value = 99.41087106270125
$ws.Range("E31").Value = 98.9146265434442
$ws.Range("E32").Value = 98.91463763602438
$ws.Range("F33").Value = 98.42379250006694
$ws.Range("F34").Value = 98.49264714748637
$ws.Range("G35").Value = 97.40635471068954
$ws.Range("G36").Value = 97.42778747851091
$ws.Range("B37").Value = 67.60636684597291
$ws.Range("B38").Value = 86.59598141447368
$ws.Range("H39").Value = 96.04247423065073
$ws.Range("H40").Value = 95.96634544343867
$ws.Range("C41").Value = 95.7121846452328
$ws.Range("C42").Value = 94.97249845532599
$ws.Range("D43").Value = 99.44703514348899
$ws.Range("D44").Value = 99.34577735510128
$ws.Range("E45").Value = 98.91138185675827
$ws.Range("E46").Value = 98.83089295312048
$ws.Range("F47").Value = 98.48285042339023
$ws.Range("F48").Value = 98.33735736014239
$ws.Range("G49").Value = 97.41540552714444
$ws.Range("G50").Value = 97.50055532518626
